$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user")

# Insert two new columns (date_of_birth, gender) right after "user_given_name" (col C),
# pushing client_since/nationality/pm_id/saving_account/trading_account two columns right.
$ws.Columns("D:E").Insert()

# New header cells
$ws.Range("D1").Value2 = "date_of_birth"
$ws.Range("E1").Value2 = "gender"

# Fix two mis-spelled surnames while we're at it
$ws.Range("B3").Value2 = "Mariano"
$ws.Range("B4").Value2 = "Ukeja"

# Fill in date_of_birth (as real dates) and gender for every user row
$ws.Range("D2").Value2 = 34832
$ws.Range("D3").Value2 = 38258
$ws.Range("D4").Value2 = 23730
$ws.Range("D5").Value2 = 30470

$ws.Range("E2").Value2 = "male"
$ws.Range("E3").Value2 = "male"
$ws.Range("E4").Value2 = "female"
$ws.Range("E5").Value2 = "male"

# Correcting client_since (now in column F) so it is consistent with the
# newly-added date_of_birth values
$ws.Range("F3").Value2 = 2022
$ws.Range("F4").Value2 = 1996

# Give the new date_of_birth column the same date format already used
# elsewhere in the workbook (reuse existing style rather than inventing one)
[void]$wb.Worksheets.Item("pending_appointments").Range("A2").Copy()
[void]$ws.Range("D2:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the column widths shown in the final layout
$ws.Columns("D").ColumnWidth = 12.1640625
$ws.Columns("E").ColumnWidth = 10.1640625

# The edit was made with the "user" sheet active, cell B4 selected
[void]$ws.Activate()
[void]$ws.Range("B4").Select()
